$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, derived from the authoritative diff.
$updates = [ordered]@{
    'D2' = '69.625.29'
    'E2' = '  +0.13%  '
    'D3' = '2.501.49'
    'E3' = '  -0.36%  '
    'E4' = '  +0.05%  '
    'D5' = '575.58'
    'E5' = '  -0.54%  '
    'D6' = '166.89'
    'E6' = '  +0.00%  '
    'E7' = '  +0.01%  '
    'E8' = '  -1.93%  '
    'D9' = '2.500.52'
    'E9' = '  -0.34%  '
    'E10' = '  +0.71%  '
    'E11' = '  +0.09%  '
    'E12' = '  +3.03%  '
    'E13' = '  +1.54%  '
    'D14' = '2.962.39'
    'D15' = '69.532.38'
    'E15' = '  +0.21%  '
    'E16' = '  +1.64%  '
    'D17' = '24.75'
    'E17' = '  -0.87%  '
    'D18' = '2.502.05'
    'E18' = '  +0.09%  '
    'E19' = '  -0.80%  '
    'D20' = '7.48'
    'E20' = '  -4.11%  '
    'D21' = '348.54'
    'E21' = '  -0.11%  '
    'E22' = '  -0.75%  '
    'E23' = '  +0.36%  '
    'E24' = '  -0.05%  '
    'D25' = '70.61'
    'E25' = '  +2.75%  '
    'E26' = '  -0.64%  '
    'D27' = '8.78'
    'E27' = '  -1.72%  '
    'D28' = '2.632.26'
    'E28' = '  -0.20%  '
    'E29' = '  +1.00%  '
    'D30' = '0.0₃0892'
    'E30' = '  -1.19%  '
    'E31' = '  -0.82%  '
    'D32' = '456.29'
    'E32' = '  -1.60%  '
    'D33' = '1.21'
    'E33' = '  -3.60%  '
    'D35' = '1.00'
    'E35' = '  +0.09%  '
    'E36' = '  -1.96%  '
    'D37' = '157.41'
    'E37' = '  +2.72%  '
    'D38' = '19.07'
    'E38' = '  +0.56%  '
    'E40' = '  +0.06%  '
    'E41' = '  -0.18%  '
    'E42' = '  -1.41%  '
    'E43' = '  -0.18%  '
    'D44' = '38.20'
    'E44' = '  +0.29%  '
    'D45' = '2.20'
    'E45' = '  -4.00%  '
    'E46' = '  -6.48%  '
    'D47' = '141.74'
    'E47' = '  -1.03%  '
    'E48' = '  -0.47%  '
    'E49' = '  -1.95%  '
    'E50' = '  +0.14%  '
    'E51' = '  -0.83%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so values like '69.625.29', '1.00', '38.20' or
    # '  +0.13%  ' are not silently re-interpreted as numbers/percentages.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
